$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENCIMENTO PRODUTOS")

# Clear out the sample/placeholder stock rows (goiaba / laranja) so the
# sheet goes back to being blank like the rows below it.
$ws.Range("B4:H4").ClearContents()
$ws.Range("J4").ClearContents()

$ws.Range("B5:H5").ClearContents()
$ws.Range("J5").ClearContents()
